$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A7").Value = "Content"
$ws.Range("B7").Value = 1.0158
$ws.Range("C7").Value = 0.7882

$ws.Range("D4").Select()
